# Improved Let error reporting.
#
# The "name =" paragraph held a single let-field (`m:v`) whose variable
# lookup failed. It is replaced by a properly formed `m:let` field (with
# a leading _GoBack bookmark, matching Word's usual behaviour) reporting
# a clearer error message, followed by a new paragraph holding the
# matching `m:endlet` closing field.

$d = $word.ActiveDocument

# Locate the paragraph to replace: the one starting with "name =" that
# contains the old `m:v` field and the red error message.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "name*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    # Fall back to the known position (2nd paragraph) if the text-based
    # lookup above ever fails to match.
    $target = $d.Paragraphs.Item(2)
}

$newParagraphsXml = @'
    <w:p>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:fldChar w:fldCharType="begin"/>
      </w:r>
      <w:r>
        <w:instrText>m:</w:instrText>
      </w:r>
      <w:r>
        <w:instrText>let</w:instrText>
      </w:r>
      <w:r>
        <w:instrText xml:space="preserve"> </w:instrText>
      </w:r>
      <w:r>
        <w:instrText xml:space="preserve"> </w:instrText>
      </w:r>
      <w:r>
        <w:instrText>=</w:instrText>
      </w:r>
      <w:r>
        <w:instrText xml:space="preserve"> self.</w:instrText>
      </w:r>
      <w:r>
        <w:instrText>name</w:instrText>
      </w:r>
      <w:r>
        <w:fldChar w:fldCharType="end"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b w:val="true"/>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>Invalid let statement: Missing identifier</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:fldChar w:fldCharType="begin"/>
      </w:r>
      <w:r>
        <w:instrText xml:space="preserve"> </w:instrText>
      </w:r>
      <w:r>
        <w:instrText>m:</w:instrText>
      </w:r>
      <w:r>
        <w:instrText>endlet</w:instrText>
      </w:r>
      <w:r>
        <w:instrText xml:space="preserve"> </w:instrText>
      </w:r>
      <w:r>
        <w:fldChar w:fldCharType="end"/>
      </w:r>
    </w:p>
'@

$packageXml = '<?xml version="1.0" standalone="yes"?>' +
    '<?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $newParagraphsXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

# Replacing the whole paragraph range with the two new paragraphs both
# removes the old field/run content and inserts the new ones in a single
# structural edit.
$target.Range.InsertXML($packageXml)
